$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content additions (bio info, images, sample data) ---
$ws.Range("E2").Value = "N"

$ws.Range("B3").Value = "Kuranda Tree Frog"
$ws.Range("E3").Value = "Y"
$ws.Range("H3").Value = "Common Name?"

$ws.Range("B4").Value = "Peron's Tree Frog"
$ws.Range("H4").Value = "Common Name?"

$ws.Range("H5").Value = "Y"

$ws.Range("H6").Value = "Y"

$ws.Range("B7").Value = "Orange Thighed Tree Frog"
$ws.Range("H7").Value = "Common Name?"

$ws.Range("H8").Value = "Y"

$ws.Range("H9").Value = "Y"

$ws.Range("H10").Value = "Y"

$ws.Range("B11").Value = "Purple-crowned fairy wren"
$ws.Range("H11").Value = "Common Name?"

# --- Column width change (col H / 8) ---
$ws.Columns.Item(8).ColumnWidth = 17.08984375

# --- View / pane changes ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B12").Select()

$excel.ActiveWindow.Left = 120
$excel.ActiveWindow.Top = 1950
$excel.ActiveWindow.Width = 25600
